$wb = $excel.ActiveWorkbook

# --- Farmland sheet: mark rows 4-9 (rows 5:11) as not included ---
$farmland = $wb.Worksheets.Item("Farmland")

$farmland.Range("D5:D11").Value = "N"

$farmland.Range("D5:D11").Select()

# --- Crops sheet: un-plan a handful of crops and let the "Plant?" ---
# --- auto-filter (shows only "Y") hide the rows that no longer match ---
$crops = $wb.Worksheets.Item("Crops")

$crops.Range("X24").Value = "N"
$crops.Range("X25").Value = "N"
$crops.Range("X26").Value = "N"
$crops.Range("X27").Value = "N"
$crops.Range("X44").Value = "N"

$crops.Rows.Item(24).Hidden = $true
$crops.Rows.Item(25).Hidden = $true
$crops.Rows.Item(26).Hidden = $true
$crops.Rows.Item(27).Hidden = $true
$crops.Rows.Item(38).Hidden = $true
$crops.Rows.Item(44).Hidden = $true

$crops.Range("X47").Select()
$crops.Activate()
